$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 167, shifting existing rows 167:239 down to 168:240
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new data point
$ws.Cells.Item(167, 1).Value = 45050
$ws.Cells.Item(167, 2).Value = 117.34

# Copy the date style (column A, style index 1 => numFmtId 14 date format)
# from the row below (now row 168, which held the previous row-167 data)
# onto the new row's date cell.
$ws.Cells.Item(168, 1).Copy() | Out-Null
$ws.Cells.Item(167, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the value (PasteSpecial formats only, but be safe and re-set it)
$ws.Cells.Item(167, 1).Value = 45050

# Update the view: scrolled to A151, with B167 selected
$ws.Activate()
$ws.Range("B167").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 151
$excel.ActiveWindow.ScrollColumn = 1

# Set the page setup (paper size 9 = A4, orientation = portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
